# Enhance the link extraction from website,
# Mark the "Status" column (D) as "Done" for rows 171 through 390
# (journal entries whose link extraction has now been processed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 171; $row -le 390; $row++) {
    $ws.Cells.Item($row, 4).Value = "Done"
}
